$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K / "Strike#" replaced values), rows 2-35
$gValues = @{
    2  = 0
    3  = 5
    4  = 3
    5  = 5
    6  = 3
    7  = 2
    8  = 4
    9  = 8
    10 = 2
    11 = 1
    12 = 3
    13 = 2
    14 = 5
    15 = 2
    16 = 3
    17 = 5
    18 = 6
    19 = 5
    20 = 4
    21 = 2
    22 = 1
    23 = 3
    24 = 5
    25 = 8
    26 = 4
    27 = 3
    28 = 1
    29 = 2
    30 = 4
    31 = 1
    32 = 4
    33 = 2
    34 = 1
    35 = 3
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
